$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has columns A (Rank), B (Team), C (ExpPoints).
# We need to insert 4 new columns (WIN, TOP2, TOP4, RELEGATION) between
# Team and ExpPoints, so ExpPoints moves from column C to column G.
$ws.Range("C:F").Insert()

# Label the newly inserted header cells.
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP2"
$ws.Range("E1").Value = "TOP4"
$ws.Range("F1").Value = "RELEGATION"

# Match the header formatting (bold, centered, bordered) used by the rest
# of row 1 by copying the format from the existing "Rank" header cell.
$ws.Range("A1").Copy()
$ws.Range("C1:F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# The new WIN/TOP2/TOP4/RELEGATION columns are placeholders for the
# upcoming Monte Carlo simulation results - leave the data rows blank for
# now, just touch them so the (empty) cells are materialised.
$ws.Range("C2:F19").Font.Bold = $false
